# Fruta / hortaliza, semanal
# Insert a new daily price record for "Vega Modelo de Temuco - Mango" right
# before the existing row 231 (which pushes the current row 231 and every
# row after it down by one row) and populate it with the new observation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift row 231 (and everything below it) down by one row.
$ws.Rows(231).Insert()

# Fill in the new row 231 with the new weekly observation.
$ws.Cells.Item(231, 1).Value  = 10
$ws.Cells.Item(231, 2).Value  = "Vega Modelo de Temuco"
$ws.Cells.Item(231, 3).Value  = "La Araucanía"
$ws.Cells.Item(231, 4).Value  = 44627
$ws.Cells.Item(231, 5).Value  = 9
$ws.Cells.Item(231, 6).Value  = "Fruta"
$ws.Cells.Item(231, 7).Value  = 100108
$ws.Cells.Item(231, 8).Value  = "Tropicales y subtropicales"
$ws.Cells.Item(231, 9).Value  = 100108002
$ws.Cells.Item(231, 10).Value = "Mango"
$ws.Cells.Item(231, 11).Value = "Sin especificar"
$ws.Cells.Item(231, 12).Value = "Primera"
$ws.Cells.Item(231, 13).Value = 500
$ws.Cells.Item(231, 14).Value = 7000
$ws.Cells.Item(231, 15).Value = 8000
$ws.Cells.Item(231, 16).Value = 7600
$ws.Cells.Item(231, 17).Value = "$/bandeja 4 kilos"
$ws.Cells.Item(231, 18).Value = "Perú"
$ws.Cells.Item(231, 19).Value = 1900
$ws.Cells.Item(231, 20).Value = 4
